# Trade #54 closed at 2026-02-18 00:21:01 - unknown UNKNOWN +0.000%
#
# This script mirrors the bot's two effects on the workbook:
#  1) An already-open HighProbConvergence trade (row 84 in "All Trades",
#     row 8 in "HighProbConvergence") is closed out -> exit price, P&L,
#     capital-after, exit reason and duration get filled in.
#  2) A brand new "momentum" trade (#112) is opened -> appended as a new
#     row at the bottom of "All Trades" (row 113) and of "momentum" (row 31).
#  3) The rollups on "Summary" and "Strategy Status" are refreshed to match.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    # Writes a literal string into a cell without letting the host's
    # "smart" input parsing reinterpret date-shaped text (e.g. "2026-02-18")
    # as a serial date number / date-formatted cell.
    param($ws, $row, $col, $text)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Summary sheet rollups
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.64
$summary.Range("B4").Value = 0.75
$summary.Range("B6").Value = 82
$summary.Range("B7").Value = 42
$summary.Range("B9").Value = 51.22

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - HighProbConvergence row (row 3)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C3").Value = 100.25
$status.Range("D3").Value = 7
$status.Range("E3").Value = 0.26
$status.Range("F3").Value = 0.25
$status.Range("G3").Value = 85.70999999999999

# ---------------------------------------------------------------------
# 3. All Trades sheet - close out trade #83 (row 84)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(84, 7).Value = 0.86        # G: Exit Price
$allTrades.Cells.Item(84, 8).Value = "CLOSED"    # H: Status
$allTrades.Cells.Item(84, 9).Value = 1.1765      # I: P&L %
$allTrades.Cells.Item(84, 10).Value = 0.01       # J: P&L $
$allTrades.Cells.Item(84, 11).Value = 100.25     # K: Capital After
Set-TextCell $allTrades 84 12 "early_exit"       # L: Exit Reason
$allTrades.Cells.Item(84, 13).Value = 0.12       # M: Duration (min)

# ---------------------------------------------------------------------
# 4. All Trades sheet - append new trade #112 (row 113)
# ---------------------------------------------------------------------
$allTrades.Cells.Item(113, 1).Value = 112
Set-TextCell $allTrades 113 2 "2026-02-18"
Set-TextCell $allTrades 113 3 "00:20:56"
Set-TextCell $allTrades 113 4 "momentum"
Set-TextCell $allTrades 113 5 "DOWN"
$allTrades.Cells.Item(113, 6).Value = 0.85
$allTrades.Cells.Item(113, 8).Value = "OPEN"
$allTrades.Cells.Item(113, 9).Value = 0
$allTrades.Cells.Item(113, 10).Value = 0
$allTrades.Cells.Item(113, 11).Value = 99.7087371310913
$allTrades.Cells.Item(113, 13).Value = 0
$allTrades.Cells.Item(113, 14).Value = 0
$allTrades.Cells.Item(113, 15).Value = 0
$allTrades.Cells.Item(113, 16).Value = 0.9
Set-TextCell $allTrades 113 17 "Downward momentum: -3.810% over 10 samples"

# ---------------------------------------------------------------------
# 5. momentum sheet - append same new trade #112 (row 31)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(31, 1).Value = 112
Set-TextCell $momentum 31 2 "2026-02-18"
Set-TextCell $momentum 31 3 "00:20:56"
Set-TextCell $momentum 31 4 "momentum"
Set-TextCell $momentum 31 5 "DOWN"
$momentum.Cells.Item(31, 6).Value = 0.85
$momentum.Cells.Item(31, 8).Value = "OPEN"
$momentum.Cells.Item(31, 9).Value = 0
$momentum.Cells.Item(31, 10).Value = 0
$momentum.Cells.Item(31, 11).Value = 99.7087371310913
$momentum.Cells.Item(31, 12).Value = 0
$momentum.Cells.Item(31, 13).Value = 0
$momentum.Cells.Item(31, 14).Value = 0.9
Set-TextCell $momentum 31 15 "Downward momentum: -3.810% over 10 samples"
$momentum.Cells.Item(31, 17).Value = 0

# ---------------------------------------------------------------------
# 6. HighProbConvergence sheet - close out trade #83 (row 8)
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(8, 7).Value = 0.86          # G: Exit Price
$hpc.Cells.Item(8, 8).Value = "CLOSED"      # H: Status
$hpc.Cells.Item(8, 9).Value = 1.1765        # I: P&L %
$hpc.Cells.Item(8, 10).Value = 0.01         # J: P&L $
$hpc.Cells.Item(8, 11).Value = 100.25       # K: Capital After
Set-TextCell $hpc 8 16 "early_exit"         # P: Exit Reason
$hpc.Cells.Item(8, 17).Value = 0.12         # Q: Duration (min)
